{"js": "// 1) Move the (empty) \"_GoBack\" bookmark out of the tail of the letterhead\n//    paragraph and into the blank paragraph that immediately follows it.\n//    This is the OOXML effect of the author placing the cursor right before\n//    the bookmark and pressing Enter: the paragraph gets split there, and\n//    the (zero-length) bookmark ends up alone in the newly-separated blank\n//    paragraph while the original paragraph keeps only its trailing tab.\nconst bookmarkName = \"_GoBack\";\nconst bookmarkRangeForLookup = context.document.getBookmarkRangeOrNullObject(bookmarkName);\nbookmarkRangeForLookup.load(\"isNullObject\");\nawait context.sync();\n\nif (!bookmarkRangeForLookup.isNullObject) {\n  const hostParagraph = bookmarkRangeForLookup.paragraphs.getFirst();\n  const destinationParagraph = hostParagraph.getNext();\n\n  // Remove the bookmark from its current (end-of-paragraph) position.\n  context.document.deleteBookmark(bookmarkName);\n  await context.sync();\n\n  // Re-create it inside the following (blank) paragraph, using the\n  // paragraph's \"Content\" range so both the start and end markers land\n  // inside that same paragraph instead of spilling into the next one.\n  const destinationContent = destinationParagraph.getRange(\"Content\");\n  destinationContent.insertBookmark(bookmarkName);\n  await context.sync();\n}\n\n// 2) Update the two VML \"textpath\" date stamps from \"Marzo 2022\" to\n//    \"Abril 2022\". Those strings live inside legacy <w:pict>/VML markup\n//    (a rotated rubber-stamp style shape), which isn't reachable through\n//    the InlinePicture/Shape object model, so we patch the paragraph's\n//    underlying OOXML directly.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  const paragraphRange = paragraph.getRange();\n  const ooxml = paragraphRange.getOoxml();\n  await context.sync();\n\n  if (ooxml.value.indexOf(\"Marzo\") !== -1) {\n    const updatedXml = ooxml.value.split(\"Marzo\").join(\"Abril\");\n    paragraphRange.insertOoxml(updatedXml, Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Move the (empty) \"_GoBack\" bookmark out of the tail of the letterhead\n#    paragraph and into the blank paragraph that immediately follows it.\n#    This mirrors the author placing the cursor right before the bookmark\n#    and pressing Enter: the paragraph splits there, so the bookmark ends\n#    up alone in the newly-separated blank paragraph while the original\n#    paragraph keeps only its trailing tab.\n$bookmarkName = \"_GoBack\"\nif ($d.Bookmarks.Exists($bookmarkName)) {\n    $bookmark = $d.Bookmarks($bookmarkName)\n    $hostParagraph = $bookmark.Range.Paragraphs(1)\n\n    # Locate the paragraph immediately following the bookmark's host\n    # paragraph using range arithmetic (no hardcoded paragraph index).\n    $afterHostPosition = $hostParagraph.Range.End\n    $destinationParagraph = $d.Range($afterHostPosition, $afterHostPosition).Paragraphs(1)\n\n    $bookmark.Delete()\n    $d.Bookmarks.Add($bookmarkName, $destinationParagraph.Range)\n}\n\n# 2) Update the two VML \"textpath\" date stamps from \"Marzo 2022\" to\n#    \"Abril 2022\". Those strings live inside legacy VML <w:pict> markup\n#    (a rotated rubber-stamp style shape) that isn't reachable through the\n#    Shapes/InlineShapes object model, so each paragraph's raw OOXML is\n#    patched directly wherever the old month name shows up.\nforeach ($paragraph in $d.Paragraphs) {\n    $paragraphXml = $paragraph.Range.WordOpenXML\n    if ($paragraphXml -like \"*Marzo*\") {\n        $updatedXml = $paragraphXml.Replace(\"Marzo\", \"Abril\")\n        $paragraph.Range.InsertXML($updatedXml)\n    }\n}\n"}
